# "Generate Report for Archive"
#
# The localization status moves from "Ready for handoff" to
# "In Translation" on every sheet (Overview shows it twice, once per
# target locale column; each locale sheet shows it once in its Status
# column). Because the new text is shorter, the Status column narrows
# on each sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn status is column E, de-de status is column F ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: Status is column C ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: Status is column C ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.5
